# Generate Report for Handback
# Updates the handback-status report timestamps for the
# "2cb4c16e-5412-4c80-892d-9095861aed28" file after a new handback run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G) for the
# 2cb4c16e-5412-4c80-892d-9095861aed28.md row (row 2).
$overview.Range("G2").Value = "2016-09-06 07:58:17"

# zh-cn sheet, row 2 (2cb4c16e-5412-4c80-892d-9095861aed28.md):
#   H = Correspond Handoff Datetime
#   K = Correspond Handback DateTime
$zhcn.Range("H2").Value = "2016-09-06 07:57:59"
$zhcn.Range("K2").Value = "2016-09-06 07:58:53"

# de-de sheet, row 2 (2cb4c16e-5412-4c80-892d-9095861aed28.md):
#   K = Correspond Handback DateTime
$dede.Range("K2").Value = "2016-09-06 07:59:17"
